$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd gender value first ("maile" -> "male"), then populate the
# full_name column which was missing a value. Order matters here only for
# the resulting shared-string table layout, matching how the values were
# entered originally.
$ws.Range("D2").Value = "male"
$ws.Range("A2").Value = "Joe Strummer"

# Restore the selection to a single cell (A2) instead of the whole row.
$ws.Range("A2").Select()
